$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, [string]$Text)
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text
    $scratch.Copy()
    $ws.Range($CellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

$ws.Range("D2").Value = "24.811.74"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.708.21"
$ws.Range("E3").Value = "  +0.48%  "
Set-TextValue "D4" "0.9956"
$ws.Range("E4").Value = "  -0.70%  "
Set-TextValue "D5" "317.88"
$ws.Range("E5").Value = "  +0.60%  "
Set-TextValue "D6" "0.9971"
$ws.Range("E6").Value = "  -0.62%  "
Set-TextValue "D7" "0.3932"
$ws.Range("E7").Value = "  +0.05%  "
Set-TextValue "D8" "0.4076"
$ws.Range("E8").Value = "  +0.68%  "
Set-TextValue "D9" "1.500"
$ws.Range("E9").Value = "  -1.13%  "
Set-TextValue "D10" "53.99"
$ws.Range("E10").Value = "  +2.12%  "
Set-TextValue "D11" "0.9955"
$ws.Range("E11").Value = "  -0.74%  "
Set-TextValue "D12" "0.08840"
$ws.Range("E12").Value = "  -0.27%  "
Set-TextValue "D13" "26.35"
$ws.Range("E13").Value = "  +11.40%  "
Set-TextValue "D14" "7.493"
$ws.Range("E14").Value = "  +0.99%  "
Set-TextValue "D15" "8.152"
$ws.Range("E15").Value = "  +0.39%  "
Set-TextValue "D16" "0.00001365"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").Value = "1.698.84"
$ws.Range("E17").Value = "  -0.47%  "
Set-TextValue "D18" "97.74"
$ws.Range("E18").Value = "  -1.67%  "
Set-TextValue "D19" "0.07172"
$ws.Range("E19").Value = "  +1.66%  "
Set-TextValue "D20" "20.62"
$ws.Range("E20").Value = "  +4.00%  "
Set-TextValue "D21" "7.320"
$ws.Range("E21").Value = "  +3.55%  "
Set-TextValue "D22" "0.9964"
$ws.Range("E22").Value = "  -0.99%  "
Set-TextValue "D23" "14.41"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").Value = "24.788.27"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("E25").Value = "  -3.70%  "
Set-TextValue "D26" "2.338"
$ws.Range("E26").Value = "  -0.53%  "
Set-TextValue "D27" "23.06"
$ws.Range("E27").Value = "  +1.80%  "
Set-TextValue "D28" "167.70"
$ws.Range("E28").Value = "  +2.09%  "
Set-TextValue "D29" "5.882"
$ws.Range("E29").Value = "  +14.07%  "
Set-TextValue "D30" "8.611"
$ws.Range("E30").Value = "  -3.62%  "
Set-TextValue "D31" "145.24"
$ws.Range("E31").Value = "  +7.15%  "
$ws.Range("D32").Value = "1.886.57"
$ws.Range("E32").Value = "  -0.39%  "
Set-TextValue "D33" "0.08849"
$ws.Range("E33").Value = "  -2.10%  "
Set-TextValue "D34" "2.181"
$ws.Range("E34").Value = "  +11.01%  "
Set-TextValue "D35" "1.076"
$ws.Range("E35").Value = "  +0.38%  "
Set-TextValue "D36" "7.227"
$ws.Range("E36").Value = "  -6.39%  "
Set-TextValue "D37" "0.03133"
$ws.Range("E37").Value = "  +5.79%  "
Set-TextValue "D38" "0.2823"
$ws.Range("E38").Value = "  +2.12%  "
Set-TextValue "D39" "0.8541"
$ws.Range("E39").Value = "  +10.39%  "
Set-TextValue "D40" "10.97"
$ws.Range("E40").Value = "  -0.44%  "
Set-TextValue "D41" "0.09229"
$ws.Range("E41").Value = "  -0.13%  "
Set-TextValue "D42" "14.23"
$ws.Range("E42").Value = "  -1.56%  "
Set-TextValue "D43" "1.479"
$ws.Range("E43").Value = "  +0.60%  "
Set-TextValue "D44" "17.66"
$ws.Range("E44").Value = "  +9.24%  "
Set-TextValue "D45" "2.733"
$ws.Range("E45").Value = "  +5.19%  "
Set-TextValue "D46" "0.7523"
$ws.Range("E46").Value = "  +4.54%  "
Set-TextValue "D47" "4.281"
$ws.Range("E47").Value = "  +1.67%  "
Set-TextValue "D48" "1.393"
$ws.Range("E48").Value = "  +2.58%  "
Set-TextValue "D49" "0.9967"
$ws.Range("E49").Value = "  -0.55%  "
Set-TextValue "D50" "140.76"
$ws.Range("E50").Value = "  +0.55%  "
Set-TextValue "D51" "0.08265"
$ws.Range("E51").Value = "  +3.60%  "
